{"js": "// Changelog edit: add a new bullet \"Ilvl of items on character panel\"\n// right after the \"BossesKilled in LFG panel \" line (same run formatting:\n// Helvetica, 12pt / sz=24).\n\nconst body = context.document.body;\n\n// Locate the anchor paragraph by its (stable, whitespace-trimmed) text.\nconst searchText = \"BossesKilled in LFG panel\";\nconst results = body.search(searchText, { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(`Anchor text not found: \"${searchText}\"`);\n}\n\n// Get the paragraph that contains the found range, then insert a new\n// paragraph right after it.\nconst anchorParagraphs = results.items[0].paragraphs;\nanchorParagraphs.load(\"text\");\nawait context.sync();\n\nconst anchorParagraph = anchorParagraphs.items[0];\nconst insertionRange = anchorParagraph.getRange(Word.RangeLocation.end);\n\n// Use insertOoxml so the new run carries xml:space=\"preserve\" exactly like\n// the rest of the document's runs, and formatting (Helvetica, sz 24 / 12pt)\n// matches the sibling paragraphs verbatim.\nconst newParagraphOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr/>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:ascii=\"Helvetica\" w:hAnsi=\"Helvetica\" w:cs=\"Helvetica\"/>\n                <w:sz w:val=\"24\"/>\n                <w:sz-cs w:val=\"24\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\">Ilvl of items on character panel</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ninsertionRange.insertOoxml(newParagraphOoxml, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Changelog edit: add a new bullet \"Ilvl of items on character panel\"\n# right after the \"BossesKilled in LFG panel \" line (same run formatting:\n# Helvetica, 12pt / sz=24).\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph by its (stable) text.\n$anchorIndex = -1\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*BossesKilled in LFG panel*\") {\n        $anchorIndex = $i\n        break\n    }\n    $i++\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Anchor paragraph 'BossesKilled in LFG panel' not found\"\n}\n\n# 1-based COM collection index for the anchor paragraph.\n$anchorPara = $d.Paragraphs.Item($anchorIndex + 1)\n$anchorRange = $anchorPara.Range\n$anchorRange.Collapse(0)            # wdCollapseEnd\n$anchorRange.InsertParagraphAfter() # creates a new empty paragraph right after,\n                                     # inheriting the anchor run's formatting\n                                     # (Helvetica / sz 24 / sz-cs 24).\n\n# Re-resolve the freshly created paragraph (immediately after the anchor) and\n# give it its text.\n$newPara = $d.Paragraphs.Item($anchorIndex + 2)\n$newPara.Range.Text = \"Ilvl of items on character panel\"\n"}
